# Automatische test-sync: 2025-08-05 19:23:50
#
# 1) Append a new logged e-mail as row 43 on the "Logs" sheet and extend
#    the sheet's conditional-formatting ranges (D/G/H/I/J columns) so they
#    keep covering the data through the new row.
# 2) Re-rank the "Overig" / "Documentatie / Datasheets" / "Kwaliteit /
#    Certificaten" rows on the "Dashboard" sheet to reflect the new counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Logs": add row 43
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A43").Value = "Kun je dit voor me fixen?"
$logs.Range("B43").Value = "mailmind.test@zohomail.eu"
$logs.Range("C43").Value = "Testmail #2: Kun je dit voor me fixen?"
$logs.Range("D43").Value = "Overig"
$logs.Range("E43").Value = "Beste [naam],
Bedankt voor je bericht. Om je beter te kunnen helpen, heb ik wat meer informatie nodig. Kun je specifiek aangeven wat er gefixt moet worden en eventueel ook wat meer details geven over het probleem dat je ervaart? Hoe meer informatie je kunt verstrekken, hoe beter we je kunnen assisteren.
Ik zie graag je reactie tegemoet.
Met vriendelijke groet,
[Naam]  
E-mailassistent"
$logs.Range("F43").Value = "2025-08-05 19:23:32"
$logs.Range("G43").Value = "Ja"
$logs.Range("H43").Value = "Nee"
$logs.Range("I43").Value = "Ja"
$logs.Range("J43").Value = "Nee"

# Extend the conditional formatting ranges from row 42 to row 43 so the
# newly added row is covered too, keeping all existing rule settings.
$ranges = @("D2:D42", "G2:G42", "H2:H42", "I2:I42", "J2:J42")
foreach ($addr in $ranges) {
    $col = $addr.Substring(0, 1)
    $newRange = $logs.Range($col + "2:" + $col + "43")
    $fc = $logs.Range($addr).FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# Sheet "Dashboard": re-sort the tail of the category breakdown
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A8").Value = "Overig"
$dash.Range("B8").Value = 2
$dash.Range("A9").Value = "Documentatie / Datasheets"
$dash.Range("B9").Value = 1
$dash.Range("A10").Value = "Kwaliteit / Certificaten"
$dash.Range("B10").Value = 1
